$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 rows belonging to the worker who was dropped from this
# statement (PEDRO VICENTE VELANDIA QUINTERO, doc 13847161). Deleting the
# entire rows shifts everything below up and keeps formatting/merged cells
# (e.g. the bottom signature block, and the special border on the last
# data row) consistent with the rest of the table.
$ws.Range("B30:J34").EntireRow.Delete() | Out-Null

# Updated detail rows (Tipo Doc, N Doc, Nombre, Periodo Mora, Valor Mora, Salario Basico)
$data = @(
    @(16, "9154751", "ADOLFO PEREZ ATENCIO", "1611", 27582, 689545),
    @(17, "73579363", "JUAN CARLOS GOMEZ PINTO", "1611", 27582, 689545),
    @(18, "1047426503", "JOSE LUIS POLO JULIO", "1611", 27582, 689545),
    @(19, "73122381", "FRANCISCO MORENO MORENO", "1611", 27582, 689545),
    @(20, "73116113", "JESUS VALENCIA FERNANDEZ", "1611", 27582, 689545),
    @(21, "1047412879", "CARLOS ANDRES MARTINEZ MARTINEZ", "1611", 27582, 689545),
    @(22, "9154751", "ADOLFO PEREZ ATENCIO", "1612", 27582, 689545),
    @(23, "73579363", "JUAN CARLOS GOMEZ PINTO", "1612", 27582, 689545),
    @(24, "1047426503", "JOSE LUIS POLO JULIO", "1612", 27582, 689545),
    @(25, "73122381", "FRANCISCO MORENO MORENO", "1612", 27582, 689545),
    @(26, "73116113", "JESUS VALENCIA FERNANDEZ", "1612", 27582, 689545),
    @(27, "1047412879", "CARLOS ANDRES MARTINEZ MARTINEZ", "1612", 27582, 689545),
    @(28, "11171600", "SILFREDO SOLERA RAMOS", "1701", 28440, 711000),
    @(29, "9154751", "ADOLFO PEREZ ATENCIO", "1701", 27582, 689545),
    @(30, "73579363", "JUAN CARLOS GOMEZ PINTO", "1701", 27582, 689545),
    @(31, "1047426503", "JOSE LUIS POLO JULIO", "1701", 27582, 689545),
    @(32, "73122381", "FRANCISCO MORENO MORENO", "1701", 27582, 689545),
    @(33, "73116113", "JESUS VALENCIA FERNANDEZ", "1701", 27582, 689545),
    @(34, "1047465978", "JAIME HENRIQUE MARIMON MARTINEZ", "1701", 28440, 711000),
    @(35, "9185722", "ALBEIRO LUIS AVILEZ DE LA ROSA", "1701", 28440, 711000),
    @(36, "9185968", "ARIEL MERCADO GUERRERO", "1701", 28440, 711000),
    @(37, "1047412879", "CARLOS ANDRES MARTINEZ MARTINEZ", "1701", 27582, 689545),
    @(38, "73168716", "RAFAEL IGNACIO AVILES DE LA ROSA", "1701", 28440, 711000),
    @(39, "9154751", "ADOLFO PEREZ ATENCIO", "1702", 27582, 689545),
    @(40, "73579363", "JUAN CARLOS GOMEZ PINTO", "1702", 27582, 689545),
    @(41, "1047426503", "JOSE LUIS POLO JULIO", "1702", 27582, 689545),
    @(42, "73122381", "FRANCISCO MORENO MORENO", "1702", 27582, 689545),
    @(43, "73116113", "JESUS VALENCIA FERNANDEZ", "1702", 27582, 689545),
    @(44, "1047465978", "JAIME HENRIQUE MARIMON MARTINEZ", "1702", 28440, 711000),
    @(45, "9185722", "ALBEIRO LUIS AVILEZ DE LA ROSA", "1702", 28440, 711000),
    @(46, "9185968", "ARIEL MERCADO GUERRERO", "1702", 28440, 711000),
    @(47, "1047412879", "CARLOS ANDRES MARTINEZ MARTINEZ", "1702", 27582, 689545),
    @(48, "73168716", "RAFAEL IGNACIO AVILES DE LA ROSA", "1702", 28440, 711000),
    @(49, "9154751", "ADOLFO PEREZ ATENCIO", "1703", 27582, 689545),
    @(50, "1047426503", "JOSE LUIS POLO JULIO", "1703", 27582, 689545),
    @(51, "73122381", "FRANCISCO MORENO MORENO", "1703", 27582, 689545),
    @(52, "73116113", "JESUS VALENCIA FERNANDEZ", "1703", 27582, 689545),
    @(53, "1047465978", "JAIME HENRIQUE MARIMON MARTINEZ", "1703", 28440, 711000),
    @(54, "9185722", "ALBEIRO LUIS AVILEZ DE LA ROSA", "1703", 28440, 711000),
    @(55, "9185968", "ARIEL MERCADO GUERRERO", "1703", 28440, 711000),
    @(56, "1047412879", "CARLOS ANDRES MARTINEZ MARTINEZ", "1703", 27582, 689545),
    @(57, "73168716", "RAFAEL IGNACIO AVILES DE LA ROSA", "1703", 28440, 711000),
    @(58, "9154751", "ADOLFO PEREZ ATENCIO", "1704", 27582, 689545),
    @(59, "73579363", "JUAN CARLOS GOMEZ PINTO", "1704", 27582, 689545),
    @(60, "1047426503", "JOSE LUIS POLO JULIO", "1704", 27582, 689545),
    @(61, "73122381", "FRANCISCO MORENO MORENO", "1704", 27582, 689545),
    @(62, "73116113", "JESUS VALENCIA FERNANDEZ", "1704", 27582, 689545),
    @(63, "1047465978", "JAIME HENRIQUE MARIMON MARTINEZ", "1704", 28440, 711000),
    @(64, "9185722", "ALBEIRO LUIS AVILEZ DE LA ROSA", "1704", 28440, 711000),
    @(65, "9185968", "ARIEL MERCADO GUERRERO", "1704", 28440, 711000),
    @(66, "1047412879", "CARLOS ANDRES MARTINEZ MARTINEZ", "1704", 27582, 689545),
    @(67, "73168716", "RAFAEL IGNACIO AVILES DE LA ROSA", "1704", 28440, 711000),
    @(68, "9154751", "ADOLFO PEREZ ATENCIO", "1705", 27582, 689545),
    @(69, "73579363", "JUAN CARLOS GOMEZ PINTO", "1705", 27582, 689545),
    @(70, "1047426503", "JOSE LUIS POLO JULIO", "1705", 27582, 689545),
    @(71, "73122381", "FRANCISCO MORENO MORENO", "1705", 27582, 689545),
    @(72, "73116113", "JESUS VALENCIA FERNANDEZ", "1705", 27582, 689545),
    @(73, "1047465978", "JAIME HENRIQUE MARIMON MARTINEZ", "1705", 28440, 711000),
    @(74, "9185722", "ALBEIRO LUIS AVILEZ DE LA ROSA", "1705", 28440, 711000),
    @(75, "9185968", "ARIEL MERCADO GUERRERO", "1705", 28440, 711000),
    @(76, "1047412879", "CARLOS ANDRES MARTINEZ MARTINEZ", "1705", 27582, 689545),
    @(77, "73168716", "RAFAEL IGNACIO AVILES DE LA ROSA", "1705", 28440, 711000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value2 = $row[1]
    $ws.Cells.Item($r, 4).Value2 = $row[2]
    $ws.Cells.Item($r, 5).Value2 = $row[3]
    $ws.Cells.Item($r, 6).Value2 = $row[4]
    $ws.Cells.Item($r, 7).Value2 = $row[5]
}

# Refresh the summary figures at the top of the statement
$ws.Range("E11").Value2 = 1728102
$ws.Range("C13").Value2 = 11
